$p = $ppt.ActivePresentation

# Slide 8 ("The Process") - title text is split into two runs: "Our " + "Process".
$s8 = $p.Slides.Item(8)
$s8.Shapes.Item(1).TextFrame.TextRange.Text = "Our Process"

# Delete slide 9 ("What It Takes")
$p.Slides.Item(9).Delete()

# Delete what is now slide 9 ("Space to Meet", originally slide 10)
$p.Slides.Item(9).Delete()
